$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the LAN block (rows 1-4): 192.168.0.x -> 201.0.0.x ---
$ws.Range("C1").Value = "201.0.0.0  - 201.0.0.7 "
$ws.Range("G1").Value = "201.0.0.1"
$ws.Range("C2").Value = "201.0.0.2"
$ws.Range("C3").Value = "201.0.0.3"
$ws.Range("C4").Value = "201.0.0.4"

# --- Update the DMZ block (rows 6-7): 192.168.1.x -> 200.0.0.x ---
$ws.Range("C6").Value = "200.0.0.0 - 200.0.0.2"
$ws.Range("G6").Value = "200.0.0.1"
$ws.Range("C7").Value = "200.0.0.2"

# --- Add the DNS block (rows 9-10) ---
$ws.Range("A9").Font.Bold = $true
$ws.Range("A9").Value = "DNS"

$ws.Range("C9").HorizontalAlignment = -4108
$ws.Range("C9").Value = "8.0.0.0 - 8.0.0.2"

$ws.Range("D9").HorizontalAlignment = -4108
$ws.Range("D9").NumberFormat = "#,##0"
$ws.Range("D9").Value = 255255255252

$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").Value = "/30"

$ws.Range("F9").HorizontalAlignment = -4152
$ws.Range("F9").Value = "Gateway:"

$ws.Range("G9").HorizontalAlignment = -4108
$ws.Range("G9").Value = "8.0.0.1"

$ws.Range("A10").Font.Italic = $true
$ws.Range("A10").Value = "DNS Server"

$ws.Range("C10").HorizontalAlignment = -4108
$ws.Range("C10").Font.Underline = $true
$ws.Range("C10").Value = "8.0.0.2"

# --- Add the Router block (rows 12-14) ---
$ws.Range("A12").Font.Bold = $true
$ws.Range("A12").Value = "Router"

$ws.Range("A13").Font.Italic = $true
$ws.Range("A13").Value = "Router InfoPoint"

$ws.Range("C13").HorizontalAlignment = -4108
$ws.Range("C13").Font.Underline = $true
$ws.Range("C13").Font.Italic = $true
$ws.Range("C13").Value = "20.0.0.1"

$ws.Range("A14").Font.Italic = $true
$ws.Range("A14").Value = "ISP"

$ws.Range("C14").HorizontalAlignment = -4108
$ws.Range("C14").Font.Underline = $true
$ws.Range("C14").Font.Italic = $true
$ws.Range("C14").Value = "20.0.0.2"

# --- Add the Strada block (rows 16-17) ---
$ws.Range("A16").Font.Bold = $true
$ws.Range("A16").Value = "Strada"

$ws.Range("C16").HorizontalAlignment = -4108
$ws.Range("C16").Value = "180.0.0.1 - 180.0.1.254"

$ws.Range("D16").HorizontalAlignment = -4108
$ws.Range("D16").Value = "255.255.254.0"

$ws.Range("E16").HorizontalAlignment = -4108
$ws.Range("E16").Value = "/23"

$ws.Range("F16").HorizontalAlignment = -4152
$ws.Range("F16").Value = "Gateway:"

$ws.Range("G16").HorizontalAlignment = -4108
$ws.Range("G16").Value = "180.0.0.0"

$ws.Range("A17").Value = "Tablet Cliente"

$ws.Range("C17").HorizontalAlignment = -4108
$ws.Range("C17").Font.Underline = $true
$ws.Range("C17").Value = "180.0.0.2"

# --- Update the selected cell (cursor position) ---
$ws.Range("F23").Select() | Out-Null
